$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new columns I ("I0") and J ("IF") ---
# Copy the formatting from the existing header cell (H1) onto the two new
# header cells so they pick up the same bold/centered/bordered style, then
# set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-28: new values for columns I and J ---
$iValues = @(5, 8, 7, 6, 9, 6, 7, 10, 8, 6, 7, 6, 7, 7, 8, 10, 5, 9, 6, 5, 5, 7, 4, 9, 6, 9, 5)
$jValues = @(6, 9, 7, 7, 9, 6, 7, 10, 8, 6, 7, 6, 7, 7, 8, 10, 6, 9, 6, 5, 5, 7, 5, 9, 6, 9, 5)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
